$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M44").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
